# Auto-generated Excel COM-interop script to apply cell value updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (83 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 288.46667
$ws.Cells.Item(2, 9).Value = 305.2143
$ws.Cells.Item(2, 11).Value = 305.2143
$ws.Cells.Item(2, 13).Value = -192.2143
$ws.Cells.Item(11, 8).Value = 35.77778
$ws.Cells.Item(11, 9).Value = 35.77778
$ws.Cells.Item(11, 11).Value = 35.77778
$ws.Cells.Item(11, 13).Value = 104.22222
$ws.Cells.Item(64, 8).Value = 7978.1035
$ws.Cells.Item(64, 9).Value = 3523.8572
$ws.Cells.Item(64, 10).Value = 9395.362999999999
$ws.Cells.Item(64, 11).Value = 3523.8572
$ws.Cells.Item(64, 12).Value = 9395.362999999999
$ws.Cells.Item(64, 13).Value = -3275.8572
$ws.Cells.Item(64, 14).Value = -9891.362999999999
$ws.Cells.Item(67, 8).Value = 7978.1035
$ws.Cells.Item(67, 9).Value = 3523.8572
$ws.Cells.Item(67, 10).Value = 9395.362999999999
$ws.Cells.Item(67, 11).Value = 3523.8572
$ws.Cells.Item(67, 12).Value = 9395.362999999999
$ws.Cells.Item(67, 13).Value = -2665.8572
$ws.Cells.Item(67, 14).Value = -11111.363
$ws.Cells.Item(80, 8).Value = 348.83334
$ws.Cells.Item(80, 9).Value = 180.2
$ws.Cells.Item(80, 10).Value = 469.2857
$ws.Cells.Item(80, 11).Value = 540.5999999999999
$ws.Cells.Item(80, 12).Value = 1407.8571
$ws.Cells.Item(80, 13).Value = 457.4000000000001
$ws.Cells.Item(80, 14).Value = -3403.8571
$ws.Cells.Item(83, 8).Value = 348.83334
$ws.Cells.Item(83, 9).Value = 180.2
$ws.Cells.Item(83, 10).Value = 469.2857
$ws.Cells.Item(83, 11).Value = 1621.8
$ws.Cells.Item(83, 12).Value = 4223.571300000001
$ws.Cells.Item(83, 13).Value = 3370.2
$ws.Cells.Item(83, 14).Value = -14207.5713
$ws.Cells.Item(86, 8).Value = 4611.4443
$ws.Cells.Item(86, 9).Value = 4666.6665
$ws.Cells.Item(86, 10).Value = 4501
$ws.Cells.Item(86, 11).Value = 4666.6665
$ws.Cells.Item(86, 12).Value = 4501
$ws.Cells.Item(86, 13).Value = -3543.6665
$ws.Cells.Item(86, 14).Value = -6747
$ws.Cells.Item(88, 8).Value = 4033.7334
$ws.Cells.Item(88, 9).Value = 687.5
$ws.Cells.Item(88, 10).Value = 4548.5386
$ws.Cells.Item(88, 11).Value = 687.5
$ws.Cells.Item(88, 12).Value = 4548.5386
$ws.Cells.Item(88, 13).Value = -281.5
$ws.Cells.Item(88, 14).Value = -5360.5386
$ws.Cells.Item(89, 8).Value = 4611.4443
$ws.Cells.Item(89, 9).Value = 4666.6665
$ws.Cells.Item(89, 10).Value = 4501
$ws.Cells.Item(89, 11).Value = 23333.3325
$ws.Cells.Item(89, 12).Value = 22505
$ws.Cells.Item(89, 13).Value = -17717.3325
$ws.Cells.Item(89, 14).Value = -33737
$ws.Cells.Item(91, 8).Value = 4033.7334
$ws.Cells.Item(91, 9).Value = 687.5
$ws.Cells.Item(91, 10).Value = 4548.5386
$ws.Cells.Item(91, 11).Value = 687.5
$ws.Cells.Item(91, 12).Value = 4548.5386
$ws.Cells.Item(91, 13).Value = 716.5
$ws.Cells.Item(91, 14).Value = -7356.5386
$ws.Cells.Item(96, 8).Value = 770182.75
$ws.Cells.Item(96, 9).Value = 1111657.5
$ws.Cells.Item(96, 11).Value = 3334972.5
$ws.Cells.Item(96, 13).Value = -3333599.5
$ws.Cells.Item(100, 8).Value = 4761.3335
$ws.Cells.Item(100, 9).Value = 1822.1428
$ws.Cells.Item(100, 11).Value = 1822.1428
$ws.Cells.Item(100, 13).Value = -1281.1428
$ws.Cells.Item(103, 8).Value = 1246.125
$ws.Cells.Item(103, 10).Value = 1298.5
$ws.Cells.Item(103, 12).Value = 3895.5
$ws.Cells.Item(103, 14).Value = -5067.5
$ws.Cells.Item(129, 8).Value = 1887.2727
$ws.Cells.Item(129, 9).Value = 822.8570999999999
$ws.Cells.Item(129, 10).Value = 3750
$ws.Cells.Item(129, 11).Value = 2468.5713
$ws.Cells.Item(129, 12).Value = 11250
$ws.Cells.Item(129, 13).Value = 2531.4287
$ws.Cells.Item(129, 14).Value = -21250

# --- Sheet: ARM (4 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2629.6365
$ws.Cells.Item(132, 9).Value = 2058.2273
$ws.Cells.Item(132, 11).Value = 6174.6819
$ws.Cells.Item(132, 13).Value = -3644.6819

# --- Sheet: BSM (22 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1615.5
$ws.Cells.Item(20, 10).Value = 1353.4286
$ws.Cells.Item(20, 12).Value = 1353.4286
$ws.Cells.Item(20, 14).Value = -1847.4286
$ws.Cells.Item(86, 8).Value = 2416
$ws.Cells.Item(86, 9).Value = 2261.4285
$ws.Cells.Item(86, 10).Value = 2570.5715
$ws.Cells.Item(86, 11).Value = 2261.4285
$ws.Cells.Item(86, 12).Value = 2570.5715
$ws.Cells.Item(86, 13).Value = -1138.4285
$ws.Cells.Item(86, 14).Value = -4816.5715
$ws.Cells.Item(89, 8).Value = 2416
$ws.Cells.Item(89, 9).Value = 2261.4285
$ws.Cells.Item(89, 10).Value = 2570.5715
$ws.Cells.Item(89, 11).Value = 11307.1425
$ws.Cells.Item(89, 12).Value = 12852.8575
$ws.Cells.Item(89, 13).Value = -5691.1425
$ws.Cells.Item(89, 14).Value = -24084.8575
$ws.Cells.Item(134, 8).Value = 2529.647
$ws.Cells.Item(134, 9).Value = 2529.647
$ws.Cells.Item(134, 11).Value = 7588.941
$ws.Cells.Item(134, 13).Value = -5053.941

# --- Sheet: CRP (4 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(97, 8).Value = 123089.5
$ws.Cells.Item(97, 10).Value = 185000
$ws.Cells.Item(97, 12).Value = 185000
$ws.Cells.Item(97, 14).Value = -186982

# --- Sheet: CUL (14 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(60, 8).Value = 47619436
$ws.Cells.Item(60, 9).Value = 55555676
$ws.Cells.Item(60, 10).Value = 2000
$ws.Cells.Item(60, 11).Value = 166667028
$ws.Cells.Item(60, 12).Value = 6000
$ws.Cells.Item(60, 13).Value = -166666777
$ws.Cells.Item(60, 14).Value = -6502
$ws.Cells.Item(113, 8).Value = 1572.3846
$ws.Cells.Item(113, 9).Value = 710.5
$ws.Cells.Item(113, 10).Value = 1955.4445
$ws.Cells.Item(113, 11).Value = 2131.5
$ws.Cells.Item(113, 12).Value = 5866.333500000001
$ws.Cells.Item(113, 13).Value = 38.5
$ws.Cells.Item(113, 14).Value = -10206.3335

# --- Sheet: GSM (30 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 22666
$ws.Cells.Item(63, 9).Value = 10000
$ws.Cells.Item(63, 10).Value = 28999
$ws.Cells.Item(63, 11).Value = 10000
$ws.Cells.Item(63, 12).Value = 28999
$ws.Cells.Item(63, 13).Value = -9314
$ws.Cells.Item(63, 14).Value = -30371
$ws.Cells.Item(66, 8).Value = 22666
$ws.Cells.Item(66, 9).Value = 10000
$ws.Cells.Item(66, 10).Value = 28999
$ws.Cells.Item(66, 11).Value = 30000
$ws.Cells.Item(66, 12).Value = 86997
$ws.Cells.Item(66, 13).Value = -26568
$ws.Cells.Item(66, 14).Value = -93861
$ws.Cells.Item(96, 8).Value = 50000.5
$ws.Cells.Item(96, 10).Value = 50000.5
$ws.Cells.Item(96, 12).Value = 50000.5
$ws.Cells.Item(96, 14).Value = -55492.5
$ws.Cells.Item(107, 8).Value = 404.73334
$ws.Cells.Item(107, 9).Value = 405.2
$ws.Cells.Item(107, 11).Value = 405.2
$ws.Cells.Item(107, 13).Value = 1514.8
$ws.Cells.Item(126, 8).Value = 2254.75
$ws.Cells.Item(126, 9).Value = 2350
$ws.Cells.Item(126, 11).Value = 7050
$ws.Cells.Item(126, 13).Value = -4580
$ws.Cells.Item(132, 8).Value = 3397.5
$ws.Cells.Item(132, 9).Value = 3524.5454
$ws.Cells.Item(132, 11).Value = 10573.6362
$ws.Cells.Item(132, 13).Value = -8043.636200000001

# --- Sheet: LTW (26 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 7081.9443
$ws.Cells.Item(40, 9).Value = 4860
$ws.Cells.Item(40, 10).Value = 8859.5
$ws.Cells.Item(40, 11).Value = 4860
$ws.Cells.Item(40, 12).Value = 8859.5
$ws.Cells.Item(40, 13).Value = -4724
$ws.Cells.Item(40, 14).Value = -9131.5
$ws.Cells.Item(68, 8).Value = 11949.1
$ws.Cells.Item(68, 9).Value = 2350.5
$ws.Cells.Item(68, 11).Value = 2350.5
$ws.Cells.Item(68, 13).Value = -1601.5
$ws.Cells.Item(71, 8).Value = 11949.1
$ws.Cells.Item(71, 9).Value = 2350.5
$ws.Cells.Item(71, 11).Value = 11752.5
$ws.Cells.Item(71, 13).Value = -8008.5
$ws.Cells.Item(117, 8).Value = 130987
$ws.Cells.Item(117, 10).Value = 130987
$ws.Cells.Item(117, 12).Value = 130987
$ws.Cells.Item(117, 14).Value = -140165
$ws.Cells.Item(136, 8).Value = 4883.0586
$ws.Cells.Item(136, 9).Value = 5021
$ws.Cells.Item(136, 10).Value = 4552
$ws.Cells.Item(136, 11).Value = 15063
$ws.Cells.Item(136, 12).Value = 13656
$ws.Cells.Item(136, 13).Value = -12513
$ws.Cells.Item(136, 14).Value = -18756

# --- Sheet: WVR (57 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 33438332
$ws.Cells.Item(5, 10).Value = 50021250
$ws.Cells.Item(5, 12).Value = 50021250
$ws.Cells.Item(5, 14).Value = -50021474
$ws.Cells.Item(19, 8).Value = 7499.5
$ws.Cells.Item(19, 10).Value = 7499.5
$ws.Cells.Item(19, 12).Value = 7499.5
$ws.Cells.Item(19, 14).Value = -7847.5
$ws.Cells.Item(62, 8).Value = 10521.429
$ws.Cells.Item(62, 10).Value = 12800
$ws.Cells.Item(62, 12).Value = 12800
$ws.Cells.Item(62, 14).Value = -14048
$ws.Cells.Item(65, 8).Value = 10521.429
$ws.Cells.Item(65, 10).Value = 12800
$ws.Cells.Item(65, 12).Value = 64000
$ws.Cells.Item(65, 14).Value = -70240
$ws.Cells.Item(75, 8).Value = 37500
$ws.Cells.Item(75, 10).Value = 50000
$ws.Cells.Item(75, 12).Value = 50000
$ws.Cells.Item(75, 14).Value = -51872
$ws.Cells.Item(78, 8).Value = 37500
$ws.Cells.Item(78, 10).Value = 50000
$ws.Cells.Item(78, 12).Value = 150000
$ws.Cells.Item(78, 14).Value = -159360
$ws.Cells.Item(96, 8).Value = 3012.375
$ws.Cells.Item(96, 9).Value = 800
$ws.Cells.Item(96, 10).Value = 3749.8333
$ws.Cells.Item(96, 11).Value = 800
$ws.Cells.Item(96, 12).Value = 3749.8333
$ws.Cells.Item(96, 13).Value = 573
$ws.Cells.Item(96, 14).Value = -6495.8333
$ws.Cells.Item(118, 8).Value = 159999.67
$ws.Cells.Item(118, 10).Value = 219999.5
$ws.Cells.Item(118, 12).Value = 219999.5
$ws.Cells.Item(118, 14).Value = -223313.5
$ws.Cells.Item(122, 8).Value = 4925.375
$ws.Cells.Item(122, 9).Value = 1880.6
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 5641.799999999999
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -3191.799999999999
$ws.Cells.Item(122, 14).Value = -34900
$ws.Cells.Item(126, 8).Value = 4850.4443
$ws.Cells.Item(126, 9).Value = 4850.4443
$ws.Cells.Item(126, 11).Value = 14551.3329
$ws.Cells.Item(126, 13).Value = -12081.3329
$ws.Cells.Item(132, 8).Value = 9071.647000000001
$ws.Cells.Item(132, 9).Value = 7691
$ws.Cells.Item(132, 10).Value = 10624.875
$ws.Cells.Item(132, 11).Value = 23073
$ws.Cells.Item(132, 12).Value = 31874.625
$ws.Cells.Item(132, 13).Value = -20543
$ws.Cells.Item(132, 14).Value = -36934.625
$ws.Cells.Item(136, 8).Value = 2437.625
$ws.Cells.Item(136, 9).Value = 2214.4285
$ws.Cells.Item(136, 11).Value = 6643.2855
$ws.Cells.Item(136, 13).Value = -4093.2855
